# Applies the "Add files via upload" commit:
#   - turns the plain-text URL in D71 into a real hyperlink (reusing the
#     same look as the other profile-link cells) which also grows that
#     row's height slightly;
#   - appends three new response rows (82-84) copied from the same
#     template formatting as row 81;
#   - leaves the cursor where the editor last left it (F89).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. D71: plain URL text -> live hyperlink (style + row height change) ---
$ws.Hyperlinks.Add($ws.Range("D71"), "https://www.beecrowd.com.br/judge/en/profile/948707") | Out-Null
# Hyperlinks.Add() stamps its own ad-hoc format; re-copy the canonical
# "Hyperlink" cell format from D70 (an existing hyperlinked profile cell)
# so D71 ends up sharing the same cell style as its siblings.
$ws.Range("D70").Copy()
$ws.Range("D71").PasteSpecial(-4122)
$ws.Rows.Item(71).RowHeight = 29.4

# --- 2. Append rows 82:84, formatted like the prior last row (81) ---
$ws.Range("A81:F81").Copy()
$ws.Range("A82:F84").PasteSpecial(-4122)
$ws.Rows.Item(82).RowHeight = 27.6
$ws.Rows.Item(83).RowHeight = 27.6
$ws.Rows.Item(84).RowHeight = 27.6

# Row 82 - new respondent
$ws.Range("A82").Value = 45384.645995370367
$ws.Range("B82").Value = "Rana Rohitashav Gehloch"
$ws.Range("C82").Value = "B23094"
$ws.Range("D82").Value = "https://www.beecrowd.com.br/judge/en/profile/949183"
$ws.Range("E82").Value = "CE"
$ws.Range("F82").Value = 0

# Row 83
$ws.Range("A83").Value = 45384.646585648145
$ws.Range("B83").Value = "Akshit Bhola"
$ws.Range("C83").Value = "B23112"
$ws.Range("D83").Value = "https://www.beecrowd.com.br/judge/en/profile/949151"
$ws.Range("E83").Value = "CSE"
$ws.Range("F83").Value = 0

# Row 84
$ws.Range("A84").Value = 45384.649305555555
$ws.Range("B84").Value = "Shubhankit Singh"
$ws.Range("C84").Value = "B23387"
$ws.Range("D84").Value = "https://www.beecrowd.com.br/judge/en/profile/948383"
$ws.Range("E84").Value = "MSE"
$ws.Range("F84").Value = 0

# --- 3. Restore the editor's final selection ---
$ws.Range("F89").Select() | Out-Null
